# Apply corrected (scaled-down / error-fixed) financial figures to rows 2-6
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1348
$ws.Range("E2").Value = 125
$ws.Range("F2").Value = 125
$ws.Range("G2").Value = 127
$ws.Range("H2").Value = 91
$ws.Range("I2").Value = 83
$ws.Range("J2").Value = 8
$ws.Range("K2").Value = 1386
$ws.Range("L2").Value = 599
$ws.Range("M2").Value = 787
$ws.Range("N2").Value = 710
$ws.Range("O2").Value = 77
$ws.Range("P2").Value = 217
$ws.Range("Q2").Value = 115
$ws.Range("R2").Value = -73
$ws.Range("S2").Value = -17
$ws.Range("T2").Value = 26
$ws.Range("U2").Value = 88
$ws.Range("V2").Value = 191
$ws.Range("W2").Value = 9.25
$ws.Range("X2").Value = 6.76
$ws.Range("Y2").Value = 12.46
$ws.Range("Z2").Value = 7.02
$ws.Range("AA2").Value = 76.09
$ws.Range("AB2").Value = 227.19
$ws.Range("AC2").Value = 192
$ws.Range("AD2").Value = 11.93
$ws.Range("AE2").Value = 1638
$ws.Range("AF2").Value = 1.4
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 43337615

# Row 3
$ws.Range("D3").Value = 932
$ws.Range("E3").Value = 38
$ws.Range("F3").Value = 38
$ws.Range("G3").Value = 48
$ws.Range("H3").Value = 89
$ws.Range("I3").Value = 64
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 1315
$ws.Range("L3").Value = 436
$ws.Range("M3").Value = 878
$ws.Range("N3").Value = 777
$ws.Range("O3").Value = 102
$ws.Range("P3").Value = 217
$ws.Range("Q3").Value = 185
$ws.Range("R3").Value = 23
$ws.Range("S3").Value = -24
$ws.Range("T3").Value = 10
$ws.Range("U3").Value = 175
$ws.Range("V3").Value = 167
$ws.Range("W3").Value = 4.07
$ws.Range("X3").Value = 9.539999999999999
$ws.Range("Y3").Value = 8.630000000000001
$ws.Range("Z3").Value = 6.59
$ws.Range("AA3").Value = 49.67
$ws.Range("AB3").Value = 257.7
$ws.Range("AC3").Value = 148
$ws.Range("AD3").Value = 15.94
$ws.Range("AE3").Value = 1792
$ws.Range("AF3").Value = 1.32
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 43337615

# Row 4
$ws.Range("D4").Value = 969
$ws.Range("E4").Value = 36
$ws.Range("F4").Value = 36
$ws.Range("G4").Value = 54
$ws.Range("H4").Value = 116
$ws.Range("I4").Value = 98
$ws.Range("J4").Value = 18
$ws.Range("K4").Value = 1469
$ws.Range("L4").Value = 488
$ws.Range("M4").Value = 981
$ws.Range("N4").Value = 895
$ws.Range("O4").Value = 79
$ws.Range("P4").Value = 217
$ws.Range("Q4").Value = 45
$ws.Range("R4").Value = -352
$ws.Range("S4").Value = 230
$ws.Range("T4").Value = 9
$ws.Range("U4").Value = 36
$ws.Range("V4").Value = 207
$ws.Range("W4").Value = 3.73
$ws.Range("X4").Value = 11.96
$ws.Range("Y4").Value = 11.75
$ws.Range("Z4").Value = 8.33
$ws.Range("AA4").Value = 49.76
$ws.Range("AB4").Value = 302.36
$ws.Range("AC4").Value = 227
$ws.Range("AD4").Value = 14.39
$ws.Range("AE4").Value = 2064
$ws.Range("AF4").Value = 1.58
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 43337615

# Row 5
$ws.Range("D5").Value = 1016
$ws.Range("E5").Value = 22
$ws.Range("F5").Value = 22
$ws.Range("G5").Value = -21
$ws.Range("H5").Value = -20
$ws.Range("I5").Value = -20
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 1527
$ws.Range("L5").Value = 561
$ws.Range("M5").Value = 966
$ws.Range("N5").Value = 881
$ws.Range("O5").Value = 79
$ws.Range("P5").Value = 217
$ws.Range("Q5").Value = 66
$ws.Range("R5").Value = -89
$ws.Range("S5").Value = 54
$ws.Range("T5").Value = 166
$ws.Range("U5").Value = -100
$ws.Range("V5").Value = 262
$ws.Range("W5").Value = 2.16
$ws.Range("X5").Value = -2.01
$ws.Range("Y5").Value = -2.2
$ws.Range("Z5").Value = -1.36
$ws.Range("AA5").Value = 58.01
$ws.Range("AB5").Value = 294.93
$ws.Range("AC5").Value = -45
$ws.Range("AD5").Value = -57.51
$ws.Range("AE5").Value = 2032
$ws.Range("AF5").Value = 1.27
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 43337615

# Row 6
$ws.Range("D6").Value = 1141
$ws.Range("E6").Value = 33
$ws.Range("F6").Value = 33
$ws.Range("G6").Value = 42
$ws.Range("H6").Value = 27
$ws.Range("I6").Value = 26
$ws.Range("K6").Value = 1518
$ws.Range("L6").Value = 531
$ws.Range("M6").Value = 987
$ws.Range("N6").Value = 901
$ws.Range("P6").Value = 217
$ws.Range("Q6").Value = -38
$ws.Range("R6").Value = -30
$ws.Range("S6").Value = 15
$ws.Range("T6").Value = 27
$ws.Range("U6").Value = -65
$ws.Range("V6").Value = 276
$ws.Range("W6").Value = 2.86
$ws.Range("X6").Value = 2.38
$ws.Range("Y6").Value = 2.93
$ws.Range("Z6").Value = 1.79
$ws.Range("AA6").Value = 53.75
$ws.Range("AB6").Value = 310.12
$ws.Range("AC6").Value = 60
$ws.Range("AD6").Value = 43.53
$ws.Range("AE6").Value = 2078
$ws.Range("AF6").Value = 1.26
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 43337615

# Row 6: AG/AH columns no longer populated for this row
$ws.Range("AG6:AH6").ClearContents()

# Rows 7-9: underlying source data no longer available -> clear D:AI, keep A-C
$ws.Range("D7:AI9").ClearContents()

